# Removed SMA for GPIO header placement.
#
# Pin 4 (P0.04) and Pin 5 (P0.05) no longer drive the FTDI/Button header
# ("Button Input" / "UART CTS (opt)") - they become plain spare GPIO header
# pins (GPIO1 / GPIO2). The UART RTS (opt) row and the FuelGauge ALERT (opt)
# row are removed entirely (no longer applicable once the SMA / fuel-gauge
# header is gone), and the old GPIO1-GPIO4 broken-out rows (pins 41-44,
# which used to live on the now-removed header) are removed too - those
# pins simply go back to the "NotUsed" pool along with pin 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite rows 8-16 in place -------------------------------------------------
# Helper-free, one-cell-at-a-time assignment (most reliable against the COM shim).
function Set-Row([int]$row, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
}

$i2cSDA = "I" + [char]0x00B2 + "C SDA"
$i2cSCL = "I" + [char]0x00B2 + "C SCL"
$dash = [string]([char]0x2014)

# Row 8: pin 4 / P0.04 becomes a bare GPIO1 header pin.
Set-Row 8 4 "P0.04" "GPIO1" "P0.04" "GPIO"

# Row 9: pin 5 / P0.05 becomes a bare GPIO2 header pin (replaces "UART CTS (opt)").
Set-Row 9 5 "P0.05" "GPIO2" "P0.04" "GPIO"

# Row 10: was "UART RTS (opt)" (pin 6) - now the UART RX row (pin 7), shifted up
# from old row 11 since the RTS row is dropped.
Set-Row 10 7 "P0.08" "UART RX" "MCU_RX" "From FTDI TX"

# Row 11: UART TX row (pin 8), shifted up from old row 12.
Set-Row 11 8 "P1.08" "UART TX" "MCU_TX" "To FTDI RX"

# Row 12: I2C SDA row (pin 22), shifted up from old row 14 (FuelGauge row dropped).
Set-Row 12 22 "P0.23" $i2cSDA "I2C_SDA" "To OLED/sensors"

# Row 13: I2C SCL row (pin 21), shifted up from old row 15.
Set-Row 13 21 "P0.22" $i2cSCL "I2C_SCL" $dash

# Row 14: Reset row (pin 16), shifted up from old row 16.
Set-Row 14 16 "P0.18" "Reset" "nRESET" "With pull-up"

# Row 15: SWDIO row (pin 26), shifted up from old row 17.
Set-Row 15 26 "SWDIO" "SWDIO" "SWDIO" "Debug interface"

# Row 16: SWDCLK row (pin 27), shifted up from old row 18.
Set-Row 16 27 "SWDCLK" "SWDCLK" "SWDCLK" "Debug interface"

# --- Drop the two now-unused rows right below the rewritten block --------------
# (old row 19 "UART RX" / old row 20 "UART TX" have already been consumed above
# via direct overwrite of rows 10-11, so here we only need to remove the two
# rows that used to hold "UART RTS (opt)" and "FuelGauge ALERT (opt)" from the
# tail so everything below shifts up by 2, landing the trailing NotUsed /
# pin-list rows on 24/25 instead of 26/27).
$ws.Rows(17).Delete()
$ws.Rows(17).Delete()

# --- Remove the stale SWDIO/SWDCLK duplicate rows + the old GPIO1-GPIO4 header rows ---
# After the two deletions above, the leftover rows that are no longer part of
# the table (old duplicate SWDIO/SWDCLK rows and the GPIO1-GPIO4 header rows)
# now sit at rows 17-20. Clear them out completely.
$ws.Range("A17:E20").ClearContents()

# --- Update the "NotUsed" pin list (pins 41-44 freed up, plus pin 6) -----------
$ws.Range("A24").Value = "NotUsed"
$ws.Range("A25").Value = "2,3,6,29,30,40,39,4,19,20,23,38,41,42,43,44"

# --- Tidy up the selection to mirror the saved state ---------------------------
$ws.Range("A12:XFD12").Select()
